# addressed julie's changes, and got it to 10 pages, submitting draft 1
#
# Adds two new summary columns (R = "Total Physical Activity Enjoyment"
# group averages/diff, S = "Total Mindfulness" group averages/diff) to the
# p2_sec_outcomes_kyle worksheet, computed from the existing data in
# columns H and N, and updates the view/selection state accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# R2/S2: averages of group A (rows 2-5 and 11-14)
$ws.Range("R2").Formula = "=AVERAGE(N2:N5)+AVERAGE(N11:N14)"
$ws.Range("S2").Formula = "=AVERAGE(H2:H5)+AVERAGE(H11:H14)"

# R3/S3: averages of group B (rows 6-10 and 15-21)
$ws.Range("R3").Formula = "=AVERAGE(N6:N10)+AVERAGE(N15:N21)"
$ws.Range("S3").Formula = "=AVERAGE(H6:H10) + AVERAGE(H15:H21)"

# R4/S4: difference between the two group averages
$ws.Range("R4").Formula = "=R2-R3"
$ws.Range("S4").Formula = "=S2-S3"

# Recalculate so the cached <v> values are up to date
$wb.Application.CalculateFullRebuild()

# Restore the view/selection state: scrolled down so row 5 is at the top,
# with O25 as the active selected cell.
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("O25").Select()
